# AIP-1541 / AIP-1546: switch the 3U/1U IDM-E cabling template over to the
# IDM+ calibration file, and extend the Cabling sheet's channel-name helper
# list (column W) up to 18 channels. Also trims the now-unused DSP feeder
# map "channel 4" entries (DSP_1_Feeder_Map_3 / DSP_2_Feeder_Map_3).

$wb = $excel.ActiveWorkbook

# --- 1. DeviceInfo sheet: point at the new IDM+ calibration file ---
$deviceInfo = $wb.Worksheets.Item("DeviceInfo")
$deviceInfo.Cells.Item(2, 2).Value = "3U_13I.cal"
$deviceInfo.Columns.Item(2).ColumnWidth = 19.5703125
$deviceInfo.Columns.Item(3).ColumnWidth = 11.28515625
$deviceInfo.Columns.Item(4).ColumnWidth = 12.28515625

# --- 2. Cabling sheet: populate helper column W with the channel-name list
#        used by the "label" (column B) data-validation dropdown, now
#        extended through Channel 18 ---
$cabling = $wb.Worksheets.Item("Cabling")
$channelNames = @(
    "Channel 1", "Channel 2", "Channel 3", "Channel 4", "Channel 5", "Channel 6",
    "Channel 7", "Channel 8", "Channel 9", "Channel 10", "Channel 11", "Channel 12",
    "Channel 15", "Channel 16", "Channel 17", "Channel 18"
)
$row = 6
foreach ($name in $channelNames) {
    $cabling.Cells.Item($row, 23).Value = $name
    $row = $row + 1
}

# --- 3. DSPFeederMap sheet: drop the unused DSP_1_Feeder_Map_3 /
#        DSP_2_Feeder_Map_3 rows (delete bottom row first so the other
#        row index doesn't shift before it's removed) ---
$dspFeederMap = $wb.Worksheets.Item("DSPFeederMap")
$dspFeederMap.Rows.Item(9).Delete()
$dspFeederMap.Rows.Item(5).Delete()

# --- Selections left by the authoring session ---
$deviceInfo.Cells.Select()
$cabling.Cells.Select()
$dspFeederMap.Range("A8:XFD8").Select()
$deviceInfo.Activate()

Write-Output "AIP-1541 / AIP-1546 cabling edits applied"
